$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add row 22 data: Date 45966, Total Count 576, Session Timeout Errors 16, Errors Requiring Analysis 560
$ws.Cells.Item(22, 1).Value = 45966
$ws.Cells.Item(22, 2).Value = 576
$ws.Cells.Item(22, 3).Value = 16
$ws.Cells.Item(22, 4).Value = 560

# Update selection to F18
$ws.Range("F18").Select()
